$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exo_ScenarioPara")

# The only data change: cell C2 value 1000 -> 300
$ws.Range("C2").Value = 300

# Update the active selection to match the saved view state (C2 instead of G16)
$ws.Range("C2").Select() | Out-Null
